$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5216948156410979
$ws.Range("C2").Value = 0.1607593003602119
$ws.Range("E2").Value = 0.1519330985054026
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.1717701907042866
$ws.Range("H2").Value = 0.3397003694199796
$ws.Range("I2").Value = 0.2329504130913342
$ws.Range("M2").Value = 0.2730669012474536
$ws.Range("N2").Value = 0.8616731018148371
$ws.Range("O2").Value = 0.9259042475361383
$ws.Range("B3").Value = 0.4552035214109935
$ws.Range("C3").Value = 0.146759093167617
$ws.Range("E3").Value = 0.1448084369238103
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.1697239829587218
$ws.Range("H3").Value = 0.3420254380100118
$ws.Range("I3").Value = 0.2367776845200087
$ws.Range("M3").Value = 0.2414966378619781
$ws.Range("N3").Value = 0.8621917839277486
$ws.Range("O3").Value = 0.9262051951714625
$ws.Range("B4").Value = 0.4142154061425742
$ws.Range("C4").Value = 0.1381076210527965
$ws.Range("E4").Value = 0.140570785456724
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.1686601756331214
$ws.Range("H4").Value = 0.3436517064196352
$ws.Range("I4").Value = 0.2393206894475561
$ws.Range("M4").Value = 0.2221344238830838
$ws.Range("N4").Value = 0.8628602373488121
$ws.Range("O4").Value = 0.9272061496430695
$ws.Range("B5").Value = 0.3974728340764671
$ws.Range("C5").Value = 0.1345684059686363
$ws.Range("E5").Value = 0.1388781052501358
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.1682749071413383
$ws.Range("H5").Value = 0.3443643650390413
$ws.Range("I5").Value = 0.2404054308960397
$ws.Range("M5").Value = 0.2142498820155225
$ws.Range("N5").Value = 0.8632208713843355
$ws.Range("O5").Value = 0.9278188555827711
$ws.Range("B6").Value = 0.3946903874862926
$ws.Range("C6").Value = 0.1339799032277114
$ws.Range("E6").Value = 0.1385990956194192
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.16821384192599
$ws.Range("H6").Value = 0.3444857173290288
$ws.Range("I6").Value = 0.2405884745944142
$ws.Range("M6").Value = 0.2129410102534948
$ws.Range("N6").Value = 0.8632860902908988
$ws.Range("O6").Value = 0.9279329524472928
$ws.Range("B7").Value = 0.4139897686340248
$ws.Range("C7").Value = 0.1380599450010749
$ws.Range("E7").Value = 0.1405478192237695
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.1686547846810882
$ws.Range("H7").Value = 0.343661115394525
$ws.Range("I7").Value = 0.2393351226209113
$ws.Range("M7").Value = 0.2220280668157315
$ws.Range("N7").Value = 0.8628647434247867
$ws.Range("O7").Value = 0.9272135841102056
$ws.Range("B8").Value = 0.4988029810166097
$ws.Range("C8").Value = 0.155943639435776
$ws.Range("E8").Value = 0.1494479553326045
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.171024579190707
$ws.Range("H8").Value = 0.3404608102026785
$ws.Range("I8").Value = 0.2342299322428278
$ws.Range("M8").Value = 0.2621769900891309
$ws.Range("N8").Value = 0.861779438924863
$ws.Range("O8").Value = 0.9258383600262192
$ws.Range("B9").Value = 0.6637893944923405
$ws.Range("C9").Value = 0.19056617134666
$ws.Range("E9").Value = 0.1679990040341721
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.1772084661778095
$ws.Range("H9").Value = 0.3357622353057934
$ws.Range("I9").Value = 0.2257546461852318
$ws.Range("M9").Value = 0.3410803134584341
$ws.Range("N9").Value = 0.8624190208188764
$ws.Range("O9").Value = 0.9296381479235407
$ws.Range("B10").Value = 0.7841435535194137
$ws.Range("C10").Value = 0.2157213276730943
$ws.Range("E10").Value = 0.1823151955183775
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.1827014668032803
$ws.Range("H10").Value = 0.3332731957388972
$ws.Range("I10").Value = 0.2204699743691663
$ws.Range("M10").Value = 0.3991560810428894
$ws.Range("N10").Value = 0.8645654852397797
$ws.Range("O10").Value = 0.9364211032266496
$ws.Range("B11").Value = 0.8386989862796668
$ws.Range("C11").Value = 0.2271019921744823
$ws.Range("E11").Value = 0.1889809482173845
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.185409386508482
$ws.Range("H11").Value = 0.3323503170089168
$ws.Range("I11").Value = 0.2182716396973952
$ws.Range("M11").Value = 0.4255997237518301
$ws.Range("N11").Value = 0.8659038905289833
$ws.Range("O11").Value = 0.940380220669482
$ws.Range("B12").Value = 0.8593286984240081
$ws.Range("C12").Value = 0.2314023539889263
$ws.Range("E12").Value = 0.191527407290863
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.1864650840993249
$ws.Range("H12").Value = 0.332030982789334
$ws.Range("I12").Value = 0.2174688694319293
$ws.Range("M12").Value = 0.4356167020113162
$ws.Range("N12").Value = 0.8664625696537058
$ws.Range("O12").Value = 0.9420055419220716
$ws.Range("B13").Value = 0.8548870423484232
$ws.Range("C13").Value = 0.2304766095323032
$ws.Range("E13").Value = 0.1909779867917294
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.1862363715025737
$ws.Range("H13").Value = 0.3320984161888134
$ws.Range("I13").Value = 0.2176404380376411
$ws.Range("M13").Value = 0.4334592210298069
$ws.Range("N13").Value = 0.8663399449672227
$ws.Range("O13").Value = 0.9416498844614978
$ws.Range("B14").Value = 0.8403967987403576
$ws.Range("C14").Value = 0.227455972349162
$ws.Range("E14").Value = 0.1891899989376284
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.185495631673092
$ws.Range("H14").Value = 0.3323234409771629
$ws.Range("I14").Value = 0.2182049995966224
$ws.Range("M14").Value = 0.4264237607254415
$ws.Range("N14").Value = 0.8659488152738675
$ws.Range("O14").Value = 0.9405114070774516
$ws.Range("B15").Value = 0.8315172586051176
$ws.Range("C15").Value = 0.2256045328169307
$ws.Range("E15").Value = 0.1880977139954894
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.1850458545837057
$ws.Range("H15").Value = 0.3324652011172873
$ws.Range("I15").Value = 0.2185546799068518
$ws.Range("M15").Value = 0.42211476762769
$ws.Range("N15").Value = 0.8657159842174309
$ws.Range("O15").Value = 0.9398304915205529
$ws.Range("B16").Value = 0.7805742063208641
$ws.Range("C16").Value = 0.2149762941280642
$ws.Range("E16").Value = 0.1818826776875113
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.1825287227432639
$ws.Range("H16").Value = 0.3333377241162836
$ws.Range("I16").Value = 0.2206177894144403
$ws.Range("M16").Value = 0.3974284060115565
$ws.Range("N16").Value = 0.8644852841360233
$ws.Range("O16").Value = 0.9361799839828251
$ws.Range("B17").Value = 0.7492716019569343
$ws.Range("C17").Value = 0.2084400134691009
$ws.Range("E17").Value = 0.178109379666644
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.181038245363915
$ws.Range("H17").Value = 0.3339266410368111
$ws.Range("I17").Value = 0.2219362030920813
$ws.Range("M17").Value = 0.3822902994134907
$ws.Range("N17").Value = 0.8638228570388549
$ws.Range("O17").Value = 0.9341645861666166
$ws.Range("B18").Value = 0.7312489635577322
$ws.Range("C18").Value = 0.2046746461558371
$ws.Range("E18").Value = 0.1759534879027171
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.1802006345677114
$ws.Range("H18").Value = 0.334285078276757
$ws.Range("I18").Value = 0.2227138734352678
$ws.Range("M18").Value = 0.3735855966161736
$ws.Range("N18").Value = 0.8634759280922566
$ws.Range("O18").Value = 0.9330875737904023
$ws.Range("B19").Value = 0.7251437205254092
$ws.Range("C19").Value = 0.2033987568273119
$ws.Range("E19").Value = 0.1752260073103216
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.1799204067261329
$ws.Range("H19").Value = 0.334409822923476
$ws.Range("I19").Value = 0.222980499796062
$ws.Range("M19").Value = 0.3706387412172631
$ws.Range("N19").Value = 0.8633643235157109
$ws.Range("O19").Value = 0.9327370176759473
$ws.Range("B20").Value = 0.7526057116123752
$ws.Range("C20").Value = 0.2091364212402027
$ws.Range("E20").Value = 0.1785095605706388
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.181194871791412
$ws.Range("H20").Value = 0.3338619099328639
$ws.Range("I20").Value = 0.2217938518495828
$ws.Range("M20").Value = 0.3839015364489384
$ws.Range("N20").Value = 0.8638898476521888
$ws.Range("O20").Value = 0.9343706187512311
$ws.Range("B21").Value = 0.8446537386895443
$ws.Range("C21").Value = 0.2283434597524661
$ws.Range("E21").Value = 0.1897145675880054
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.1857123819267343
$ws.Range("H21").Value = 0.3322565275246916
$ws.Range("I21").Value = 0.2180383674638655
$ws.Range("M21").Value = 0.4284901568213826
$ws.Range("N21").Value = 0.8660622937857312
$ws.Range("O21").Value = 0.9408423796531054
$ws.Range("B22").Value = 0.9046413707030183
$ws.Range("C22").Value = 0.2408423424349166
$ws.Range("E22").Value = 0.1971676851706263
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.1888413457025848
$ws.Range("H22").Value = 0.3313830091363457
$ws.Range("I22").Value = 0.2157570576827155
$ws.Range("M22").Value = 0.4576507976543098
$ws.Range("N22").Value = 0.8677842700142975
$ws.Range("O22").Value = 0.945807214882052
$ws.Range("B23").Value = 0.8726409276771392
$ws.Range("C23").Value = 0.2341764812066458
$ws.Range("E23").Value = 0.1931778398161086
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.1871551441009274
$ws.Range("H23").Value = 0.3318331370524135
$ws.Range("I23").Value = 0.2169587580934191
$ws.Range("M23").Value = 0.4420855075928074
$ws.Range("N23").Value = 0.8668376353285083
$ws.Range("O23").Value = 0.9430899638791459
$ws.Range("B24").Value = 0.7510984426272671
$ws.Range("C24").Value = 0.2088215988553941
$ws.Range("E24").Value = 0.1783285969787229
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.1811240009288824
$ws.Range("H24").Value = 0.3338911130007034
$ws.Range("I24").Value = 0.2218581475008854
$ws.Range("M24").Value = 0.3831731011677277
$ws.Range("N24").Value = 0.8638594555550299
$ws.Range("O24").Value = 0.9342772170698623
$ws.Range("B25").Value = 0.6193039122467212
$ws.Range("C25").Value = 0.1812487146434592
$ws.Range("E25").Value = 0.1628611134718341
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.1753697089420072
$ws.Range("H25").Value = 0.3368643184179518
$ws.Range("I25").Value = 0.2278824033030316
$ws.Range("M25").Value = 0.3197164188841768
$ws.Range("N25").Value = 0.8619508116596251
$ws.Range("O25").Value = 0.9279115104626072
